$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style_D2 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.711.31'
$ws.Range("D2").Style = $style_D2
$ws.Range("E2").Value = '  -2.50%  '
$style_D3 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.652.66'
$ws.Range("D3").Style = $style_D3
$ws.Range("E3").Value = '  +3.07%  '
$ws.Range("E4").Value = '  -0.32%  '
$style_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '406.09'
$ws.Range("D5").Style = $style_D5
$ws.Range("E5").Value = '  -1.42%  '
$style_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.72'
$ws.Range("D6").Style = $style_D6
$ws.Range("E6").Value = '  +3.78%  '
$style_D7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.645.97'
$ws.Range("D7").Style = $style_D7
$ws.Range("E7").Value = '  +3.05%  '
$style_D8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.619'
$ws.Range("D8").Style = $style_D8
$ws.Range("E8").Value = '  -3.01%  '
$ws.Range("E9").Value = '  +0.04%  '
$style_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.726'
$ws.Range("D10").Style = $style_D10
$ws.Range("E10").Value = '  -4.62%  '
$ws.Range("E11").Value = '  -3.48%  '
$style_D12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000324'
$ws.Range("D12").Style = $style_D12
$ws.Range("E12").Value = '  +8.62%  '
$style_D13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.80'
$ws.Range("D13").Style = $style_D13
$ws.Range("E13").Value = '  +0.56%  '
$style_D14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.90'
$ws.Range("D14").Style = $style_D14
$ws.Range("E14").Value = '  +1.92%  '
$style_D15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.255.92'
$ws.Range("D15").Style = $style_D15
$ws.Range("E15").Value = '  +3.54%  '
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$style_D17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.652.97'
$ws.Range("D17").Style = $style_D17
$ws.Range("E17").Value = '  +3.37%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$style_D18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.58'
$ws.Range("D18").Style = $style_D18
$ws.Range("E18").Value = '  +12.07%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$style_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '19.97'
$ws.Range("D19").Style = $style_D19
$ws.Range("E19").Value = '  +0.24%  '
$style_D20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.08'
$ws.Range("D20").Style = $style_D20
$ws.Range("E20").Value = '  -1.25%  '
$style_D21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '64.944.52'
$ws.Range("D21").Style = $style_D21
$ws.Range("E21").Value = '  -2.26%  '
$style_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '420.64'
$ws.Range("D22").Style = $style_D22
$ws.Range("E22").Value = '  -4.17%  '
$style_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.24'
$ws.Range("D23").Style = $style_D23
$ws.Range("E23").Value = '  +18.88%  '
$style_D24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.76'
$ws.Range("D24").Style = $style_D24
$ws.Range("E24").Value = '  -2.64%  '
$style_D25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.99'
$ws.Range("D25").Style = $style_D25
$ws.Range("E25").Value = '  -3.18%  '
$style_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '35.74'
$ws.Range("D26").Style = $style_D26
$ws.Range("E26").Value = '  +5.54%  '
$style_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.20'
$ws.Range("D27").Style = $style_D27
$ws.Range("E27").Value = '  -5.24%  '
$ws.Range("E28").Value = '  -4.89%  '
$ws.Range("E29").Value = '  +5.74%  '
$style_D30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.69'
$ws.Range("D30").Style = $style_D30
$ws.Range("E30").Value = '  +4.52%  '
$style_D31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.72'
$ws.Range("D31").Style = $style_D31
$ws.Range("E31").Value = '  -1.14%  '
$style_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.117'
$ws.Range("D32").Style = $style_D32
$ws.Range("E32").Value = '  +2.39%  '
$ws.Range("B33").Value = 'RenderToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$style_D33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.95'
$ws.Range("D33").Style = $style_D33
$ws.Range("E33").Value = '  -3.23%  '
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$style_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '41.42'
$ws.Range("D34").Style = $style_D34
$ws.Range("E34").Value = '  +6.63%  '
$ws.Range("E35").Value = '  +1.74%  '
$style_D36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.89'
$ws.Range("D36").Style = $style_D36
$ws.Range("E36").Value = '  -0.41%  '
$style_D37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = $style_D37
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("E38").Value = '  -3.45%  '
$style_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.95'
$ws.Range("D39").Style = $style_D39
$ws.Range("E39").Value = '  +32.34%  '
$ws.Range("E40").Value = '  -3.28%  '
$style_D41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.997'
$ws.Range("D41").Style = $style_D41
$ws.Range("E41").Value = '  -0.31%  '
$style_D42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0₃0653'
$ws.Range("D42").Style = $style_D42
$ws.Range("E42").Value = '  -6.16%  '
$style_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.33'
$ws.Range("D43").Style = $style_D43
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$style_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.42'
$ws.Range("D44").Style = $style_D44
$ws.Range("E44").Value = '  +4.26%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$style_D45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '26.58'
$ws.Range("D45").Style = $style_D45
$ws.Range("E45").Value = '  +27.34%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$style_D46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.12'
$ws.Range("D46").Style = $style_D46
$ws.Range("E46").Value = '  +22.95%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$style_D47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.08'
$ws.Range("D47").Style = $style_D47
$ws.Range("E47").Value = '  +7.91%  '
$style_D48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '143.79'
$ws.Range("D48").Style = $style_D48
$ws.Range("E48").Value = '  -1.63%  '
$ws.Range("E49").Value = '  -4.23%  '
$style_D50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.53'
$ws.Range("D50").Style = $style_D50
$ws.Range("E50").Value = '  -6.38%  '
$style_D51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.290'
$ws.Range("D51").Style = $style_D51
$ws.Range("E51").Value = '  -3.89%  '
